$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "1" to "Ninotsminda"
$ws.Name = "Ninotsminda"

# The left table (Total/Urban/Rural) has been "upgraded until Javakheti":
# the 2018-2020 (columns J:L) figures for the Urban and Rural rows are no
# longer reported and become confidential/unavailable placeholders ("…"),
# matching the rest of the row. The Total row (row 5) keeps its numbers.
$ws.Range("J6:L7").Value = "…"

# Remove the blank spacer row above the footnote, shifting the note up
# from row 9 to row 8 and shrinking the used range accordingly.
$ws.Rows("8").Delete()
